$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 9, column B should get the date-only format currently on B8
# (copy it first, before B8's format changes below)
$ws.Range("B8").Copy()
$ws.Range("B9").PasteSpecial(-4122)

# Row 8, column B (date 2020-06-07) changes from date-only format to datetime format
# (matching the format already used by B2:B7)
$ws.Range("B2").Copy()
$ws.Range("B8").PasteSpecial(-4122)

# New row 9, column A gets the same style as the rest of column A
$ws.Range("A8").Copy()
$ws.Range("A9").PasteSpecial(-4122)

# Fill in new row 9 values (data for June 8th, date serial 43990)
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = 43990
$ws.Range("C9").Value = 120102
$ws.Range("D9").Value = 177875
$ws.Range("E9").Value = 46398
$ws.Range("F9").Value = 14053
$ws.Range("G9").Value = 33.32
